$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update 리로's stats, remove emoji decorations from name
$ws.Range("A2").Value = "리로`n@lee-lo-4u"
$ws.Range("B2").Value = 659
$ws.Range("C2").Value = 519
$ws.Range("F2").Value = 6

# Row 3: was 김소윤 placeholder row -> becomes 조한준 with real stats
$ws.Range("A3").Value = "조한준`n@Desde_Seúl"
$ws.Range("B3").Value = 255
$ws.Range("C3").Value = 161
$ws.Range("D3").Value = 76
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 5

# Row 4: becomes 김소윤 with fire emoji, all zero stats
$ws.Range("A4").Value = "김소윤 🔥`n@catmocotto"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 5: becomes 임동한 with fire emoji, all zero stats
$ws.Range("A5").Value = "임동한 🔥`n@easyfood369"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Row 6: becomes 김수정 with fire emoji, all zero stats
$ws.Range("A6").Value = "김수정 🔥`n@kokonyang-p6l"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# Row 7: becomes 강민성, no emoji, numeric zero stats
$ws.Range("A7").Value = "강민성`n@IQ160건강깡패"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

# Row 8: becomes 김예림, no emoji, numeric zero stats
$ws.Range("A8").Value = "김예림`n@vitaminute4u"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
